# Added periodic & upfront related scenarios
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the "repaymentstrategy" value (row 17, column B) from "RBI (India)"
# to "Overdue/Due Fee/Int,Principal"
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the new active selection on that worksheet
$ws.Activate()
$ws.Range("B17").Select()
